# "Excercise-NotView: Progress in the lofic of the sequence"
#
# The sequence data shifts up by one row (the old row 1 - 159/36 - drops
# out, so old row 2 becomes new row 1, etc., leaving 10 data rows instead
# of 11), and a new (empty, underlined) helper column C is marked next to
# the first seven rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the first row of the sequence; everything below shifts up one row.
$ws.Rows(1).Delete()

# Mark column C for rows 1-7 with an underlined font (cells stay empty).
for ($r = 1; $r -le 7; $r++) {
    $ws.Cells.Item($r, 3).Font.Underline = 1
}

# Update the view: zoomed in to 175% with the selection on E4.
$excel.ActiveWindow.Zoom = 175
$ws.Range("E4").Select()
